$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Data" to "Summary"
$ws.Name = "Summary"

# Register a new named cell style "title_" (bold + underline), matching the
# updated shared style catalog. It is added to the workbook's style gallery
# but, like the original edit, is not applied to any cell on this sheet.
$newStyle = $wb.Styles.Add("title_")
$newStyle.Font.Bold = $true
$newStyle.Font.Underline = $true

# Remove row 5 (the Micro / SMEs / MSMEs header cells), which also drops the
# now-unused shared strings and shrinks the sheet's used range to A1:A3.
$ws.Rows("5:5").Delete()
